$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that were dropped from the table. Deleting from the
# bottom up keeps the row numbers of the rows still to be processed valid:
#  - rows 14 & 13: Hampton (Merrymen Cafe) x2
#  - row 11: Docklands
#  - row 10: Collingwood
#  - row 8: Cheltenham
#  - row 3: Albert Park
#  - row 2: Abbotsford
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(2).Delete()

# After the 7 deletions above, the table has shrunk from 30 data rows to
# 23 data rows, occupying rows 2-24 (header still in row 1). Append the
# two new Springvale rows at the end, in rows 25 and 26.
$ws.Cells.Item(25, 1).Value = "Springvale"
$ws.Cells.Item(25, 2).Value = "IKEA Springvale, 917 Princes Hwy"
$ws.Cells.Item(25, 3).Value = "30/12/20, 4:00pm-6:30pm"
$ws.Cells.Item(25, 4).Value = "Case shopped at store and dined at cafe"

$ws.Cells.Item(26, 1).Value = "Springvale"
$ws.Cells.Item(26, 2).Value = "Springvale Shopping Centre,  46-58 Buckingham Avenue"
$ws.Cells.Item(26, 3).Value = "29/12/20, 11:00am-12:30pm"
$ws.Cells.Item(26, 4).Value = "Case shopped"
